$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sorted data (descending by value); replaces old rows 2-23 with new rows 2-21
$data = @(
    @("English", 21.85593598916267),
    @("Chinese", 17.54876257003382),
    @("Spanish", 6.958571358246909),
    @("German", 4.470407427296483),
    @("Arabic", 4.313229381631241),
    @("Japanese", 4.288212177656546),
    @("Russian", 3.261992891270676),
    @("Malay-Indonesian", 3.042216684903867),
    @("Portuguese", 2.955798296086481),
    @("French", 2.676755852262131),
    @("Italian", 2.066725480531678),
    @("Turkish", 1.864752758099962),
    @("Korean", 1.72061643848208),
    @("Dutch", 1.253717416877074),
    @("Persian", 1.107988669990283),
    @("Thai", 0.9852287868200295),
    @("Polish", 0.9324736078057533),
    @("Urdu", 0.8648203688524005),
    @("Vietnamese", 0.6953604521218848),
    @("Bengali", 0.681861996036907)
)

# Stash a copy of the language-name cell formatting (column A, row 2) in a scratch
# cell far outside the data range so it survives the upcoming full clear.
$scratch = $ws.Range("D1")
$ws.Range("A2").Copy()
$scratch.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Wipe the entire old data range (rows 2-23), formats included, since it shrinks
# to rows 2-21 and the two trailing rows must disappear completely.
$ws.Range("A2:B23").Clear()

$rowIndex = 2
foreach ($entry in $data) {
    $nameCell = $ws.Cells.Item($rowIndex, 1)
    $scratch.Copy()
    $nameCell.PasteSpecial(-4122) # xlPasteFormats
    $nameCell.Value = $entry[0]

    $valueCell = $ws.Cells.Item($rowIndex, 2)
    $valueCell.Value = $entry[1]

    $rowIndex++
}

$scratch.Clear()
$excel.CutCopyMode = 0
